# Slide 12 ("Grammar Rules Relevant to Records"), Content Placeholder 2,
# paragraph 10:
#   "variable   = ( varId | paramId) { indexExpr | fieldExpr } ."
# becomes
#   "variable   = ( varId | paramId ) { indexExpr | fieldExpr } ."
# i.e. a space is inserted before the closing ")" so the run that used to
# read ") { " is split into " ) " and "{ ".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$para = $tr.Paragraphs(10, 1)

# Locate "paramId)" inside this paragraph and split right after the "Id",
# on the ") " that immediately follows it, inserting a leading space.
$start = $para.Start + $para.Text.IndexOf("paramId)") + ("paramId").Length

$closeParen = $tr.Characters($start, 2)
$closeParen.Text = " ) "
